$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-11 22:29:57"
$wsZhCn.Range("H3").Value = "2016-03-11 22:30:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-11 22:30:00"
$wsDeDe.Range("H3").Value = "2016-03-11 22:30:21"
